$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "Done" status for the first two requirement rows (C5, C6)
$ws.Range("C5").Value = "Done"
$ws.Range("C6").Value = "Done"

# Update the active selection to C7 to match the saved view state
$ws.Range("C7").Select()
